$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 28505.041
$ws.Range("J69").Value = 29720.143
$ws.Range("L69").Value = 89160.429
$ws.Range("N69").Value = -90908.429
$ws.Range("H70").Value = 2470.5676
$ws.Range("J70").Value = 2788.2632
$ws.Range("L70").Value = 8364.7896
$ws.Range("N70").Value = -8904.7896
$ws.Range("H72").Value = 28505.041
$ws.Range("J72").Value = 29720.143
$ws.Range("L72").Value = 267481.287
$ws.Range("N72").Value = -276217.287
$ws.Range("H73").Value = 2470.5676
$ws.Range("J73").Value = 2788.2632
$ws.Range("L73").Value = 8364.7896
$ws.Range("N73").Value = -10236.7896
$ws.Range("H107").Value = 556.6667
$ws.Range("I107").Value = 646.6667
$ws.Range("J107").Value = 466.66666
$ws.Range("K107").Value = 646.6667
$ws.Range("L107").Value = 466.66666
$ws.Range("M107").Value = 1273.3333
$ws.Range("N107").Value = -4306.66666
$ws.Range("H111").Value = 7624.875
$ws.Range("I111").Value = 3500
$ws.Range("K111").Value = 10500
$ws.Range("M111").Value = -7433
$ws.Range("H125").Value = 1556.2858
$ws.Range("I125").Value = 1588.6666
$ws.Range("K125").Value = 14297.9994
$ws.Range("M125").Value = -11837.9994
$ws.Range("H127").Value = 8380.869000000001
$ws.Range("J127").Value = 14953.272
$ws.Range("L127").Value = 44859.81600000001
$ws.Range("N127").Value = -54779.81600000001
$ws.Range("H132").Value = 11719.043
$ws.Range("I132").Value = 1100.9487
$ws.Range("J132").Value = 63482.25
$ws.Range("K132").Value = 3302.8461
$ws.Range("L132").Value = 190446.75
$ws.Range("M132").Value = -772.8460999999998
$ws.Range("N132").Value = -195506.75
$ws.Range("H137").Value = 2941.1
$ws.Range("J137").Value = 3690
$ws.Range("L137").Value = 11070
$ws.Range("N137").Value = -16170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 180.57143
$ws.Range("I5").Value = 163.66667
$ws.Range("J5").Value = 193.25
$ws.Range("K5").Value = 163.66667
$ws.Range("L5").Value = 193.25
$ws.Range("M5").Value = -51.66667000000001
$ws.Range("N5").Value = -417.25
$ws.Range("H32").Value = 9084.718000000001
$ws.Range("I32").Value = 8190.3887
$ws.Range("K32").Value = 8190.3887
$ws.Range("M32").Value = -7903.3887
$ws.Range("H45").Value = 4916.6665
$ws.Range("I45").Value = 3962.2856
$ws.Range("K45").Value = 3962.2856
$ws.Range("M45").Value = -3585.2856
$ws.Range("H74").Value = 1764.8572
$ws.Range("I74").Value = 1809.0555
$ws.Range("K74").Value = 1809.0555
$ws.Range("M74").Value = -935.0554999999999
$ws.Range("H77").Value = 1764.8572
$ws.Range("I77").Value = 1809.0555
$ws.Range("K77").Value = 9045.2775
$ws.Range("M77").Value = -4677.2775
$ws.Range("H97").Value = 1712.1428
$ws.Range("I97").Value = 1712.1428
$ws.Range("K97").Value = 1712.1428
$ws.Range("M97").Value = -1216.1428
$ws.Range("H102").Value = 1678.5714
$ws.Range("I102").Value = 1625
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1625
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -3
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 180.57143
$ws.Range("I4").Value = 163.66667
$ws.Range("J4").Value = 193.25
$ws.Range("K4").Value = 163.66667
$ws.Range("L4").Value = 193.25
$ws.Range("M4").Value = -48.66667000000001
$ws.Range("N4").Value = -423.25
$ws.Range("H20").Value = 1678.25
$ws.Range("I20").Value = 1730.75
$ws.Range("K20").Value = 1730.75
$ws.Range("M20").Value = -1483.75
$ws.Range("H99").Value = 25390.766
$ws.Range("I99").Value = 30025.715
$ws.Range("K99").Value = 30025.715
$ws.Range("M99").Value = -28527.715
$ws.Range("H134").Value = 2207.2307
$ws.Range("I134").Value = 1082.4634
$ws.Range("K134").Value = 3247.3902
$ws.Range("M134").Value = -712.3902000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2638.3777
$ws.Range("I31").Value = 1071.7778
$ws.Range("K31").Value = 1071.7778
$ws.Range("M31").Value = -776.7778000000001
$ws.Range("H34").Value = 2638.3777
$ws.Range("I34").Value = 1071.7778
$ws.Range("K34").Value = 1071.7778
$ws.Range("M34").Value = -869.7778000000001
$ws.Range("H58").Value = 1419.2812
$ws.Range("I58").Value = 1115
$ws.Range("K58").Value = 1115
$ws.Range("M58").Value = -912
$ws.Range("H86").Value = 49470.43
$ws.Range("I86").Value = 65258.8
$ws.Range("J86").Value = 9999.5
$ws.Range("K86").Value = 65258.8
$ws.Range("L86").Value = 9999.5
$ws.Range("M86").Value = -64135.8
$ws.Range("N86").Value = -12245.5
$ws.Range("H89").Value = 49470.43
$ws.Range("I89").Value = 65258.8
$ws.Range("J89").Value = 9999.5
$ws.Range("K89").Value = 326294
$ws.Range("L89").Value = 49997.5
$ws.Range("M89").Value = -320678
$ws.Range("N89").Value = -61229.5
$ws.Range("H134").Value = 2417.3555
$ws.Range("I134").Value = 1382.9722
$ws.Range("J134").Value = 6554.8887
$ws.Range("K134").Value = 4148.9166
$ws.Range("L134").Value = 19664.6661
$ws.Range("M134").Value = -1613.9166
$ws.Range("N134").Value = -24734.6661
$ws.Range("H136").Value = 1419.2812
$ws.Range("I136").Value = 1115
$ws.Range("K136").Value = 3345
$ws.Range("M136").Value = -795

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 54611390
$ws.Range("I4").Value = 69676350
$ws.Range("J4").Value = 923.875
$ws.Range("K4").Value = 209029050
$ws.Range("L4").Value = 2771.625
$ws.Range("M4").Value = -209028938
$ws.Range("N4").Value = -2995.625
$ws.Range("H7").Value = 1346.7858
$ws.Range("I7").Value = 1496.0834
$ws.Range("K7").Value = 4488.2502
$ws.Range("M7").Value = -4376.2502
$ws.Range("H26").Value = 350
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 3000
$ws.Range("N26").Value = -3576
$ws.Range("H107").Value = 300.04544
$ws.Range("J107").Value = 281.5
$ws.Range("L107").Value = 844.5
$ws.Range("N107").Value = -4684.5
$ws.Range("H136").Value = 2448.818
$ws.Range("I136").Value = 963
$ws.Range("J136").Value = 3934.6365
$ws.Range("K136").Value = 2889
$ws.Range("L136").Value = 11803.9095
$ws.Range("M136").Value = 2211
$ws.Range("N136").Value = -22003.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5000178.5
$ws.Range("I2").Value = 87.40000000000001
$ws.Range("J2").Value = 10000269
$ws.Range("K2").Value = 87.40000000000001
$ws.Range("L2").Value = 10000269
$ws.Range("M2").Value = 25.59999999999999
$ws.Range("N2").Value = -10000495
$ws.Range("H113").Value = 3792.8572
$ws.Range("I113").Value = 2726.7368
$ws.Range("J113").Value = 6043.5557
$ws.Range("K113").Value = 2726.7368
$ws.Range("L113").Value = 6043.5557
$ws.Range("M113").Value = -556.7368000000001
$ws.Range("N113").Value = -10383.5557
$ws.Range("H132").Value = 3039.5833
$ws.Range("I132").Value = 3317.6843
$ws.Range("K132").Value = 9953.052899999999
$ws.Range("M132").Value = -7423.052899999999
$ws.Range("H135").Value = 64722.223
$ws.Range("J135").Value = 64722.223
$ws.Range("L135").Value = 64722.223
$ws.Range("N135").Value = -74862.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3520.3635
$ws.Range("I7").Value = 1894.8182
$ws.Range("K7").Value = 1894.8182
$ws.Range("M7").Value = -1782.8182
$ws.Range("H22").Value = 1760
$ws.Range("I22").Value = 1266.6666
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 1266.6666
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -971.6666
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 1760
$ws.Range("I27").Value = 1266.6666
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 1266.6666
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -1159.6666
$ws.Range("N27").Value = -2714
$ws.Range("H68").Value = 5867.3706
$ws.Range("I68").Value = 4442.4443
$ws.Range("J68").Value = 6579.8335
$ws.Range("K68").Value = 4442.4443
$ws.Range("L68").Value = 6579.8335
$ws.Range("M68").Value = -3693.4443
$ws.Range("N68").Value = -8077.8335
$ws.Range("H71").Value = 5867.3706
$ws.Range("I71").Value = 4442.4443
$ws.Range("J71").Value = 6579.8335
$ws.Range("K71").Value = 22212.2215
$ws.Range("L71").Value = 32899.1675
$ws.Range("M71").Value = -18468.2215
$ws.Range("N71").Value = -40387.1675
$ws.Range("H126").Value = 3520.3635
$ws.Range("I126").Value = 1894.8182
$ws.Range("K126").Value = 5684.4546
$ws.Range("M126").Value = -3214.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2924.818
$ws.Range("I136").Value = 1682.909
$ws.Range("J136").Value = 5408.636
$ws.Range("K136").Value = 5048.727000000001
$ws.Range("L136").Value = 16225.908
$ws.Range("M136").Value = -2498.727000000001
$ws.Range("N136").Value = -21325.908
